$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.345.44'
$ws.Range("E2").Value = '  -3.10%  '

$ws.Range("D3").Value = '3.300.94'
$ws.Range("E3").Value = '  -3.76%  '

$ws.Range("E4").Value = '  +0.06%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.75'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -3.57%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.92'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -7.21%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.300.51'
$ws.Range("E8").Value = '  -3.79%  '

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.468'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -3.15%  '

$ws.Range("E10").Value = '  -2.85%  '

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.119'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -4.39%  '

$ws.Range("E12").Value = '  -2.16%  '

$ws.Range("D13").Value = '3.866.21'
$ws.Range("E13").Value = '  -3.88%  '

$ws.Range("E14").Value = '  +0.28%  '

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.87'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -6.45%  '

$ws.Range("D16").Value = '3.306.44'
$ws.Range("E16").Value = '  -3.70%  '

$ws.Range("E17").Value = '  -3.73%  '

$ws.Range("D18").Value = '60.365.31'
$ws.Range("E18").Value = '  -3.01%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.12'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  -5.35%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.02'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -3.22%  '

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.64'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -3.75%  '

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '374.86'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -2.15%  '

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.34'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -1.23%  '

$ws.Range("E24").Value = '  -0.01%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.535'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -6.21%  '

$ws.Range("D26").Value = '3.434.90'
$ws.Range("E26").Value = '  -3.45%  '

$ws.Range("E27").Value = '  -8.83%  '

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.172'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -4.63%  '

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.994'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -0.49%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.18'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -7.23%  '

$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("E32").Value = '  -4.53%  '

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.55'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  -4.97%  '

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.61'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -2.69%  '

$ws.Range("E35").Value = '  -6.85%  '

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.11'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -6.69%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.41'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("E38").Value = '  -5.07%  '

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.69'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -3.48%  '

$ws.Range("D40").Value = '3.331.80'
$ws.Range("E40").Value = '  -3.78%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.62'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -14.76%  '

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0729'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -6.62%  '

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.94'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -2.05%  '

$ws.Range("E44").Value = '  -3.70%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.13'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -6.19%  '

$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -5.27%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.57'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  -6.46%  '

$ws.Range("D48").Value = '2.347.22'
$ws.Range("E48").Value = '  -7.73%  '

$ws.Range("E49").Value = '  +0.00%  '

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.39'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -7.36%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.29'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -5.81%  '
